$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.951.27'
$ws.Range("E2").Value = '  -2.17%  '

$ws.Range("D3").Value = '1.796.94'
$ws.Range("E3").Value = '  -0.28%  '

$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").Value = "'316.94"
$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").Value = "'0.9987"
$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("D7").Value = "'0.5319"
$ws.Range("E7").Value = '  -1.00%  '

$ws.Range("D8").Value = "'0.3876"

$ws.Range("D9").Value = "'0.07451"
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").Value = "'41.41"
$ws.Range("E10").Value = '  -2.23%  '

$ws.Range("D11").Value = "'1.088"
$ws.Range("E11").Value = '  -2.26%  '

$ws.Range("D12").Value = "'0.9992"
$ws.Range("E12").Value = '  -0.24%  '

$ws.Range("D13").Value = "'6.183"
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").Value = "'7.445"
$ws.Range("E14").Value = '  +0.97%  '

$ws.Range("D15").Value = "'20.38"
$ws.Range("E15").Value = '  -1.22%  '

$ws.Range("D16").Value = '1.791.78'
$ws.Range("E16").Value = '  -0.34%  '

$ws.Range("D17").Value = "'88.43"
$ws.Range("E17").Value = '  -1.83%  '

$ws.Range("E18").Value = '  -0.46%  '

$ws.Range("D19").Value = "'0.06541"
$ws.Range("E19").Value = '  +1.54%  '

$ws.Range("D20").Value = "'0.9992"
$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("E21").Value = '  +0.13%  '

$ws.Range("D22").Value = "'5.964"
$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("D23").Value = '27.961.35'
$ws.Range("E23").Value = '  -2.26%  '

$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("D25").Value = "'2.090"
$ws.Range("E25").Value = '  -0.45%  '

$ws.Range("D26").Value = "'156.65"
$ws.Range("E26").Value = '  -1.02%  '

$ws.Range("E27").Value = '  -1.29%  '

$ws.Range("D28").Value = '1.998.40'
$ws.Range("E28").Value = '  -0.40%  '

$ws.Range("D29").Value = "'2.307"
$ws.Range("E29").Value = '  -1.87%  '

$ws.Range("D30").Value = "'121.88"
$ws.Range("E30").Value = '  -0.80%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'1.101"
$ws.Range("E31").Value = '  -0.20%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = "'0.1085"
$ws.Range("E32").Value = '  +3.07%  '

$ws.Range("D33").Value = "'3.663"
$ws.Range("E33").Value = '  -0.53%  '

$ws.Range("D34").Value = "'5.510"
$ws.Range("E34").Value = '  -2.27%  '

$ws.Range("D35").Value = "'0.06972"
$ws.Range("E35").Value = '  +7.32%  '

$ws.Range("E36").Value = '  -1.84%  '

$ws.Range("D37").Value = "'0.02273"
$ws.Range("E37").Value = '  -1.08%  '

$ws.Range("D38").Value = "'5.079"
$ws.Range("E38").Value = '  +0.95%  '

$ws.Range("E39").Value = '  -3.69%  '

$ws.Range("D40").Value = "'11.27"
$ws.Range("E40").Value = '  +0.42%  '

$ws.Range("D41").Value = "'1.193"
$ws.Range("E41").Value = '  -0.18%  '

$ws.Range("D42").Value = "'0.6117"
$ws.Range("E42").Value = '  -1.56%  '

$ws.Range("E43").Value = '  -0.65%  '

$ws.Range("D44").Value = "'13.38"
$ws.Range("E44").Value = '  +0.97%  '

$ws.Range("D45").Value = "'3.670"
$ws.Range("E45").Value = '  -0.42%  '

$ws.Range("D46").Value = "'0.5716"
$ws.Range("E46").Value = '  -2.10%  '

$ws.Range("D47").Value = "'124.52"
$ws.Range("E47").Value = '  -1.42%  '

$ws.Range("D48").Value = "'1.914"
$ws.Range("E48").Value = '  -1.24%  '

$ws.Range("E49").Value = '  +1.95%  '

$ws.Range("D50").Value = "'0.06802"
$ws.Range("E50").Value = '  -1.19%  '

$ws.Range("D51").Value = "'0.00000000295"
$ws.Range("E51").Value = '  +36.43%  '
